$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = 43967
$ws.Range("B66").Value = 42236
$ws.Range("C66").Value = 1610
$ws.Range("D66").Value = 41
$ws.Range("E66").Value = 2004

$tbl = $ws.ListObjects.Item("Table3")
$tbl.Resize($ws.Range("A1:E66"))
